# Adds two new match rows (35 and 36) to the Kuwait Premier League 2023-2024
# sheet, as produced by the scraper run on 05-11-2023 08:45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        Row = 35
        Indice = 34
        Data = 45234.64583333334
        Home = "Al Salmiya"
        HomeGols = 3
        Away = "Khaitan"
        AwayGols = 1
        HomeOpen = 1.71
        HomeOpenData = "03/11/2023 08:42"
        HomeClose = 1.47
        HomeCloseData = "04/11/2023 01:03"
        DrawOpen = 3.65
        DrawOpenData = "03/11/2023 08:42"
        DrawClose = 4.33
        DrawCloseData = "04/11/2023 13:33"
        AwayOpen = 3.69
        AwayOpenData = "03/11/2023 08:42"
        AwayClose = 5.57
        AwayCloseData = "04/11/2023 01:03"
        Url = "https://www.betexplorer.com/football/kuwait/premier-league/al-salmiya-khaitan/rgJ1JRua/"
    },
    @{
        Row = 36
        Indice = 35
        Data = 45234.76388888889
        Home = "Kazma SC"
        HomeGols = 0
        Away = "Al Qadisiya"
        AwayGols = 2
        HomeOpen = 2.64
        HomeOpenData = "03/11/2023 08:42"
        HomeClose = 3.02
        HomeCloseData = "04/11/2023 09:40"
        DrawOpen = 3.32
        DrawOpenData = "03/11/2023 08:42"
        DrawClose = 3.49
        DrawCloseData = "04/11/2023 16:26"
        AwayOpen = 2.23
        AwayOpenData = "03/11/2023 08:42"
        AwayClose = 2.12
        AwayCloseData = "04/11/2023 09:40"
        Url = "https://www.betexplorer.com/football/kuwait/premier-league/kazma-sc-al-qadisiya/vLCAH5AB/"
    }
)

# The existing data rows carry two special cell styles: column A (bold,
# bordered, centered index) and column E (date/time number format). Reuse
# them exactly by copying the formatting from the last existing row (34)
# instead of rebuilding the format from scratch.
$ws.Range("A34").Copy() | Out-Null
$ws.Range("A35:A36").PasteSpecial(-4122) | Out-Null

$ws.Range("E34").Copy() | Out-Null
$ws.Range("E35:E36").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Indice

    $ws.Cells.Item($row, 2).Value = "kuwait"
    $ws.Cells.Item($row, 3).Value = "premier-league"
    $ws.Cells.Item($row, 4).Value = "2023-2024"

    $ws.Cells.Item($row, 5).Value = $r.Data

    $ws.Cells.Item($row, 6).Value = $r.Home
    $ws.Cells.Item($row, 7).Value = $r.HomeGols
    $ws.Cells.Item($row, 8).Value = $r.Away
    $ws.Cells.Item($row, 9).Value = $r.AwayGols

    $ws.Cells.Item($row, 10).Value = $r.HomeOpen
    $ws.Cells.Item($row, 11).Value = $r.HomeOpenData
    $ws.Cells.Item($row, 12).Value = $r.HomeClose
    $ws.Cells.Item($row, 13).Value = $r.HomeCloseData

    $ws.Cells.Item($row, 14).Value = $r.DrawOpen
    $ws.Cells.Item($row, 15).Value = $r.DrawOpenData
    $ws.Cells.Item($row, 16).Value = $r.DrawClose
    $ws.Cells.Item($row, 17).Value = $r.DrawCloseData

    $ws.Cells.Item($row, 18).Value = $r.AwayOpen
    $ws.Cells.Item($row, 19).Value = $r.AwayOpenData
    $ws.Cells.Item($row, 20).Value = $r.AwayClose
    $ws.Cells.Item($row, 21).Value = $r.AwayCloseData

    $ws.Cells.Item($row, 22).Value = $r.Url
}
